# Java Key Words update
# 1) The paragraph "- definiować metody abstrakcyj|ne, ..." previously had a
#    mid-word bookmark ("_GoBack") splitting it into two runs. Remove the
#    bookmark and re-merge the text into a single run.
$d = $word.ActiveDocument

$d.Bookmarks("_GoBack").Delete()

$d.Content.Find.Execute('- definiować metody abstrakcyjne, które klasa dziedzicząca będzie musiała zaimplementować,', $true, $false, $false, $false, $false, `
                         $true, 1, $false, '- definiować metody abstrakcyjne, które klasa dziedzicząca będzie musiała zaimplementować,', 2) | Out-Null

# 2) Replace the trailing empty paragraph with the new "synchronized",
#    "volatile", "transient" and "static" keyword sections. The last
#    paragraph regains the document's "_GoBack" bookmark at its end.
$lastPara = $d.Paragraphs.Last
$newContentRange = $lastPara.Range
$newContentXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:lang w:val="en-GB"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:lang w:val="en-GB"/>
    </w:rPr>
    <w:t>synchronized</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:lang w:val="en-GB"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:lang w:val="en-GB"/>
    </w:rPr>
    <w:t xml:space="preserve">It is used to prevent concurrency. </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:val="en-GB"/>
    </w:rPr>
    <w:t>Synchronized keyword can be applied to static/non-static methods or a block of code. Only one thread at a time can access synchronized methods</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:lang w:val="en-GB"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:lang w:val="en-GB"/>
    </w:rPr>
    <w:t>volatile</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:lang w:val="en-GB"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:lang w:val="en-GB"/>
    </w:rPr>
    <w:t xml:space="preserve">It is used to mark a Java variable as “being stored in main memory”. </w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:lang w:val="en-GB"/>
    </w:rPr>
    <w:t>So</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:lang w:val="en-GB"/>
    </w:rPr>
    <w:t xml:space="preserve"> declaring a variable as </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:i/>
      <w:lang w:val="en-GB"/>
    </w:rPr>
    <w:t>volatile</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:val="en-GB"/>
    </w:rPr>
    <w:t xml:space="preserve"> thus guarantees the visibility for other threads of writes to that variable.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:lang w:val="en-GB"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:lang w:val="en-GB"/>
    </w:rPr>
    <w:t xml:space="preserve">JVM </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:lang w:val="en-GB"/>
    </w:rPr>
    <w:t>can not</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:lang w:val="en-GB"/>
    </w:rPr>
    <w:t xml:space="preserve"> reorder read and write operations of volatile variable.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:lang w:val="en-GB"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:lang w:val="en-GB"/>
    </w:rPr>
    <w:t>transient</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:lang w:val="en-GB"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:lang w:val="en-GB"/>
    </w:rPr>
    <w:t>Not required to be serialized.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:lang w:val="en-GB"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:lang w:val="en-GB"/>
    </w:rPr>
    <w:t xml:space="preserve">static </w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:lang w:val="en-GB"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:lang w:val="en-GB"/>
    </w:rPr>
    <w:t>Everything defined as static is related to class and not to object.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:lang w:val="en-GB"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:lang w:val="en-GB"/>
    </w:rPr>
    <w:t>A static method can be accessed without creating the objects. It can only access static variables.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:lang w:val="en-GB"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:lang w:val="en-GB"/>
    </w:rPr>
    <w:t>A class cannot be declared static except inner class. But a class can be said to be static if all the variables and methods of the class are static, and the constructor is private.</w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
'
$newContentRange.InsertXML($newContentXml)
